# agrego todo lo encargado
# - inserta una nueva columna "foto" en B (desplaza todo a la derecha)
# - agrega las columnas "adeuda_materias" y "quien_aprobo" al final
# - agrega las filas 3 y 4 con nuevos registros
# - actualiza "curso" (ahora columna M) de la fila 2 a "matematica"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insertar columna "foto" en B (todo lo que estaba desde B en adelante se
#    corre una columna a la derecha: DNI pasa de B a C, apellido de C a D, ...)
# ---------------------------------------------------------------------------
$ws.Columns("B:B").Insert()

$ws.Range("B1").Value = "foto"
$ws.Range("B2").Value = "https://cdn.outsideonline.com/wp-content/uploads/2023/03/Funny_Dog_H.jpg?crop=16:9&width=960&enable=upscale&quality=100"

# ---------------------------------------------------------------------------
# 2) La columna "curso" (ahora en M) de la fila 2 cambia de valor
# ---------------------------------------------------------------------------
$ws.Range("M2").Value = "matematica"

# ---------------------------------------------------------------------------
# 3) Nuevas columnas al final de la fila de encabezados
# ---------------------------------------------------------------------------
$ws.Range("AD1").Value = "adeuda_materias"
$ws.Range("AE1").Value = "quien_aprobo"

$ws.Range("AD2").Value = "ingles, matematica"
$ws.Range("AE2").Value = "susana"

# ---------------------------------------------------------------------------
# Helper: escribe un valor en una celda, replicando el formato "General"
# (estilo s="1") que ya usan las celdas existentes de la hoja, para que las
# celdas nuevas queden con el mismo formato que el resto de la tabla.
# ---------------------------------------------------------------------------
$generalFormat = $ws.Range("A1").NumberFormat

function Set-Cell($addr, $value) {
    $ws.Range($addr).NumberFormat = $generalFormat
    $ws.Range($addr).Value = $value
}

# Algunas columnas (telefono_tutor2, cuit_tutor, edad, cod_postal) guardan
# valores que parecen numeros pero deben quedar como texto (igual que en las
# filas existentes): se formatea la celda como texto antes de escribir el
# valor, y luego se restaura el formato general -- tal como haria un usuario
# que tipea esos datos en una columna formateada como texto.
function Set-TextCell($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $ws.Range($addr).NumberFormat = $generalFormat
}

# ---------------------------------------------------------------------------
# 4) Fila 3: nuevo registro (sin foto, y sin las columnas del final)
# ---------------------------------------------------------------------------
Set-Cell "A3" 2
Set-Cell "C3" "perro"
Set-Cell "D3" "apellido"
Set-Cell "E3" "perro"
Set-Cell "F3" "iud f"
Set-Cell "G3" "No"
Set-Cell "H3" "3i4343"
Set-Cell "I3" "vfoi.ju,den"
Set-Cell "J3" "o dfin"
Set-Cell "K3" "doinf "
Set-TextCell "L3" "444444444444444"
Set-Cell "M3" "programacion"
Set-Cell "N3" "este colegio"
Set-Cell "O3" "urrrr4"
Set-TextCell "P3" "341093484"
Set-Cell "Q3" "No"
Set-Cell "R3" "ninguna"
Set-Cell "S3" "No"
Set-Cell "T3" "nignuna"
Set-Cell "U3" "roro@gmail.com"
Set-Cell "V3" "noviembre"
Set-TextCell "W3" "69"
Set-Cell "X3" "san martin"
Set-Cell "Y3" "argentina"
Set-Cell "Z3" "constan"
Set-Cell "AA3" "eaviucn"
Set-TextCell "AB3" "1669"
Set-Cell "AC3" "No"

# ---------------------------------------------------------------------------
# 5) Fila 4: nuevo registro (sin foto, y sin las columnas del final)
# ---------------------------------------------------------------------------
Set-Cell "A4" 3
Set-Cell "C4" 123
Set-Cell "D4" "cocina"
Set-Cell "E4" "gato"
Set-Cell "F4" "iud f"
Set-Cell "G4" "No"
Set-Cell "H4" "3i4343"
Set-Cell "I4" "vfoi.ju,den"
Set-Cell "J4" "o dfin"
Set-Cell "K4" "doinf "
Set-TextCell "L4" "444444444444444"
Set-Cell "M4" "matematica"
Set-Cell "N4" "este colegio"
Set-Cell "O4" "urrrr4"
Set-TextCell "P4" "341093484"
Set-Cell "Q4" "No"
Set-Cell "R4" "ninguna"
Set-Cell "S4" "No"
Set-Cell "T4" "nignuna"
Set-Cell "U4" "roro@gmail.com"
Set-Cell "V4" "noviembre"
Set-TextCell "W4" "69"
Set-Cell "X4" "san martin"
Set-Cell "Y4" "argentina"
Set-Cell "Z4" "constan"
Set-Cell "AA4" "eaviucn"
Set-TextCell "AB4" "1669"
Set-Cell "AC4" "No"
